$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 881.1429000000001
$ws.Range("I38").Value = 881.1429000000001
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 2643.4287
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("M38").Value = -2271.4287
$ws.Range("H58").Value = 2716.2
$ws.Range("I58").Value = 432.4
$ws.Range("K58").Value = 1297.2
$ws.Range("M58").Value = -1147.2
$ws.Range("H98").Value = 3898.7144
$ws.Range("I98").Value = 4136.3335
$ws.Range("J98").Value = 3471
$ws.Range("K98").Value = 4136.3335
$ws.Range("L98").Value = 3471
$ws.Range("M98").Value = -2638.3335
$ws.Range("N98").Value = -6467
$ws.Range("H112").Value = 3424.3
$ws.Range("J112").Value = 3424.3
$ws.Range("L112").Value = 10272.9
$ws.Range("N112").Value = -12488.9
$ws.Range("H121").Value = 1321.625
$ws.Range("J121").Value = 1367.5714
$ws.Range("L121").Value = 4102.7142
$ws.Range("N121").Value = -7596.7142
$ws.Range("H122").Value = 3898.7144
$ws.Range("I122").Value = 4136.3335
$ws.Range("J122").Value = 3471
$ws.Range("K122").Value = 12409.0005
$ws.Range("L122").Value = 10413
$ws.Range("M122").Value = -9959.000499999998
$ws.Range("N122").Value = -15313
$ws.Range("H127").Value = 1896.375
$ws.Range("I127").Value = 1896.375
$ws.Range("K127").Value = 5689.125
$ws.Range("M127").Value = -729.125
$ws.Range("H131").Value = 3679.7778
$ws.Range("J131").Value = 4002.875
$ws.Range("L131").Value = 12008.625
$ws.Range("N131").Value = -22088.625
$ws.Range("H132").Value = 1124.1177
$ws.Range("I132").Value = 1140.6666
$ws.Range("K132").Value = 3421.9998
$ws.Range("M132").Value = -891.9998000000001
$ws.Range("H137").Value = 2604.25
$ws.Range("J137").Value = 2815.2222
$ws.Range("L137").Value = 8445.6666
$ws.Range("N137").Value = -13545.6666
$ws.Range("H138").Value = 3140.1614
$ws.Range("I138").Value = 6273.5
$ws.Range("J138").Value = 2050.3044
$ws.Range("K138").Value = 18820.5
$ws.Range("L138").Value = 6150.9132
$ws.Range("M138").Value = -13680.5
$ws.Range("N138").Value = -16430.9132
$ws.Range("H141").Value = 1275665.4
$ws.Range("I141").Value = 2002617.2
$ws.Range("K141").Value = 6007851.6
$ws.Range("M141").Value = -6002671.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3818.8965
$ws.Range("I32").Value = 3298.1135
$ws.Range("K32").Value = 3298.1135
$ws.Range("M32").Value = -3011.1135
$ws.Range("H61").Value = 2598.7334
$ws.Range("I61").Value = 1361.909
$ws.Range("K61").Value = 1361.909
$ws.Range("M61").Value = -1149.909
$ws.Range("H74").Value = 793.8333
$ws.Range("I74").Value = 673.4
$ws.Range("J74").Value = 1396
$ws.Range("K74").Value = 673.4
$ws.Range("L74").Value = 1396
$ws.Range("M74").Value = 200.6
$ws.Range("N74").Value = -3144
$ws.Range("H77").Value = 793.8333
$ws.Range("I77").Value = 673.4
$ws.Range("J77").Value = 1396
$ws.Range("K77").Value = 3367
$ws.Range("L77").Value = 6980
$ws.Range("M77").Value = 1001
$ws.Range("N77").Value = -15716
$ws.Range("H136").Value = 2598.7334
$ws.Range("I136").Value = 1361.909
$ws.Range("K136").Value = 4085.727
$ws.Range("M136").Value = -1535.727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1939.6666
$ws.Range("I20").Value = 1808.1818
$ws.Range("J20").Value = 2301.25
$ws.Range("K20").Value = 1808.1818
$ws.Range("L20").Value = 2301.25
$ws.Range("M20").Value = -1561.1818
$ws.Range("N20").Value = -2795.25
$ws.Range("H86").Value = 112615.78
$ws.Range("I86").Value = 1442.75
$ws.Range("K86").Value = 1442.75
$ws.Range("M86").Value = -319.75
$ws.Range("H89").Value = 112615.78
$ws.Range("I89").Value = 1442.75
$ws.Range("K89").Value = 7213.75
$ws.Range("M89").Value = -1597.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2253.2778
$ws.Range("I132").Value = 1734.3846
$ws.Range("J132").Value = 3602.4
$ws.Range("K132").Value = 5203.1538
$ws.Range("L132").Value = 10807.2
$ws.Range("M132").Value = -2673.1538
$ws.Range("N132").Value = -15867.2
$ws.Range("H134").Value = 1804.6316
$ws.Range("I134").Value = 1202.7142
$ws.Range("K134").Value = 3608.1426
$ws.Range("M134").Value = -1073.1426

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2072.8428
$ws.Range("J68").Value = 2126.8657
$ws.Range("L68").Value = 6380.597099999999
$ws.Range("N68").Value = -8002.597099999999
$ws.Range("H71").Value = 2072.8428
$ws.Range("J71").Value = 2126.8657
$ws.Range("L71").Value = 19141.7913
$ws.Range("N71").Value = -27253.7913
$ws.Range("H107").Value = 1480.9
$ws.Range("J107").Value = 1785.1538
$ws.Range("L107").Value = 5355.4614
$ws.Range("N107").Value = -9195.4614
$ws.Range("H137").Value = 4306.95
$ws.Range("J137").Value = 5504.615
$ws.Range("L137").Value = 16513.845
$ws.Range("N137").Value = -26713.845
$ws.Range("H140").Value = 3503.6
$ws.Range("I140").Value = 866.6429000000001
$ws.Range("K140").Value = 2599.9287
$ws.Range("M140").Value = 2580.0713
$ws.Range("H141").Value = 1989.6154
$ws.Range("I141").Value = 1989.6154
$ws.Range("K141").Value = 5968.8462
$ws.Range("M141").Value = -788.8462

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1360.3334
$ws.Range("I97").Value = 920.8823
$ws.Range("K97").Value = 920.8823
$ws.Range("M97").Value = -424.8823
$ws.Range("H122").Value = 3809.6667
$ws.Range("I122").Value = 1699.5
$ws.Range("J122").Value = 4864.75
$ws.Range("K122").Value = 5098.5
$ws.Range("L122").Value = 14594.25
$ws.Range("M122").Value = -2648.5
$ws.Range("N122").Value = -19494.25
$ws.Range("H126").Value = 1738849
$ws.Range("I126").Value = 2418097.5
$ws.Range("J126").Value = 2992
$ws.Range("K126").Value = 7254292.5
$ws.Range("L126").Value = 8976
$ws.Range("M126").Value = -7251822.5
$ws.Range("N126").Value = -13916

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4700.25
$ws.Range("J7").Value = 6499.8335
$ws.Range("L7").Value = 6499.8335
$ws.Range("N7").Value = -6723.8335
$ws.Range("H126").Value = 4700.25
$ws.Range("J126").Value = 6499.8335
$ws.Range("L126").Value = 19499.5005
$ws.Range("N126").Value = -24439.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1516
$ws.Range("I132").Value = 714.1429000000001
$ws.Range("K132").Value = 2142.4287
$ws.Range("M132").Value = 387.5712999999996
